# "researching inline nav for About page"
# Insert a new "Time spent" column into Table1 between "Time Cost (in hours)"
# and "Notes", fill in the new data for the "Inline nav" research-examples
# row, and update the sheet view / selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a blank column at D, shifting the existing "Notes" column (and all
# its data) one column to the right, then grow the table to cover it.
$ws.Range("D1:D31").Insert(-4161)
$lo.Resize($ws.Range("A1:E31"))

# Headers
$ws.Range("D1").Value = "Time spent"
$ws.Range("E1").Value = "Notes"
$ws.Range("E1").WrapText = $true

# New "time spent" note for the "Research examples" task (row 3): replace
# what used to be the "Notes" text ("Scrolling nav") with the new note, and
# record the time spent researching in the new column.
$ws.Range("D3").Value = "25 minutes"
$ws.Range("E3").Value = "Research mobile applications; determine PRL-appropriate structure; maybe ask Federica for opinion"
$ws.Range("E3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 43.5

# The other two "Notes" cells that shifted from column D to E keep their
# wrapped-text styling.
$ws.Range("E12").WrapText = $true
$ws.Range("E25").WrapText = $true

# Column widths for the new layout
$ws.Columns.Item(4).ColumnWidth = 19.33
$ws.Columns.Item(5).ColumnWidth = 31.5

# Scroll back to the top and select C3 (matches the reviewed state)
$ws.Range("C3").Select()
